$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDOG")

$ws.Range("B7").Value = 74000000.0
$ws.Range("C7").Value = 68000000.0
$ws.Range("D7").Value = 53611000.0
$ws.Range("E7").Value = 25602000.0
$ws.Range("F7").Value = 13010000.0
